$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 499131.7  # H17
$ws.Cells.Item(17, 10).Value = 499131.7  # J17
$ws.Cells.Item(17, 12).Value = 1497395.1  # L17
$ws.Cells.Item(17, 14).Value = -1497731.1  # N17
$ws.Cells.Item(33, 8).Value = 235.08696  # H33
$ws.Cells.Item(33, 9).Value = 185.38889  # I33
$ws.Cells.Item(33, 10).Value = 414  # J33
$ws.Cells.Item(33, 11).Value = 185.38889  # K33
$ws.Cells.Item(33, 12).Value = 414  # L33
$ws.Cells.Item(33, 13).Value = 43.61111  # M33
$ws.Cells.Item(33, 14).Value = -872  # N33
$ws.Cells.Item(51, 8).Value = 4670.923  # H51
$ws.Cells.Item(51, 9).Value = 0  # I51
$ws.Cells.Item(51, 10).Value = 4670.923  # J51
$ws.Cells.Item(51, 11).Value = 0  # K51
$ws.Cells.Item(51, 12).Value = 4670.923  # L51
$ws.Cells.Item(51, 13).ClearContents()  # M51
$ws.Cells.Item(51, 14).Value = -5638.923  # N51
$ws.Cells.Item(116, 8).Value = 1697.5385  # H116
$ws.Cells.Item(116, 9).Value = 1659.4286  # I116
$ws.Cells.Item(116, 10).Value = 1742  # J116
$ws.Cells.Item(116, 11).Value = 1659.4286  # K116
$ws.Cells.Item(116, 12).Value = 1742  # L116
$ws.Cells.Item(116, 13).Value = 1782.5714  # M116
$ws.Cells.Item(116, 14).Value = -8626  # N116
$ws.Cells.Item(137, 8).Value = 3031747  # H137
$ws.Cells.Item(137, 9).Value = 5883635.5  # I137
$ws.Cells.Item(137, 10).Value = 1615.625  # J137
$ws.Cells.Item(137, 11).Value = 17650906.5  # K137
$ws.Cells.Item(137, 12).Value = 4846.875  # L137
$ws.Cells.Item(137, 13).Value = -17648356.5  # M137
$ws.Cells.Item(137, 14).Value = -9946.875  # N137

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15444.816  # H32
$ws.Cells.Item(32, 9).Value = 15858.716  # I32
$ws.Cells.Item(32, 11).Value = 15858.716  # K32
$ws.Cells.Item(32, 13).Value = -15571.716  # M32
$ws.Cells.Item(74, 8).Value = 8690956  # H74
$ws.Cells.Item(74, 9).Value = 13945630  # I74
$ws.Cells.Item(74, 11).Value = 13945630  # K74
$ws.Cells.Item(74, 13).Value = -13944756  # M74
$ws.Cells.Item(77, 8).Value = 8690956  # H77
$ws.Cells.Item(77, 9).Value = 13945630  # I77
$ws.Cells.Item(77, 11).Value = 69728150  # K77
$ws.Cells.Item(77, 13).Value = -69723782  # M77
$ws.Cells.Item(80, 8).Value = 34996.668  # H80
$ws.Cells.Item(80, 10).Value = 34995  # J80
$ws.Cells.Item(80, 12).Value = 34995  # L80
$ws.Cells.Item(80, 14).Value = -36991  # N80
$ws.Cells.Item(83, 8).Value = 34996.668  # H83
$ws.Cells.Item(83, 10).Value = 34995  # J83
$ws.Cells.Item(83, 12).Value = 104985  # L83
$ws.Cells.Item(83, 14).Value = -114969  # N83
$ws.Cells.Item(102, 8).Value = 15874294  # H102
$ws.Cells.Item(102, 9).Value = 20409348  # I102
$ws.Cells.Item(102, 10).Value = 1605.5  # J102
$ws.Cells.Item(102, 11).Value = 20409348  # K102
$ws.Cells.Item(102, 12).Value = 1605.5  # L102
$ws.Cells.Item(102, 13).Value = -20407726  # M102
$ws.Cells.Item(102, 14).Value = -4849.5  # N102
$ws.Cells.Item(132, 8).Value = 67806.78  # H132
$ws.Cells.Item(132, 9).Value = 101901.9  # I132
$ws.Cells.Item(132, 10).Value = 52309  # J132
$ws.Cells.Item(132, 11).Value = 305705.7  # K132
$ws.Cells.Item(132, 12).Value = 156927  # L132
$ws.Cells.Item(132, 13).Value = -303175.7  # M132
$ws.Cells.Item(132, 14).Value = -161987  # N132

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 42000  # H59
$ws.Cells.Item(59, 10).Value = 42000  # J59
$ws.Cells.Item(59, 12).Value = 42000  # L59
$ws.Cells.Item(59, 14).Value = -43694  # N59
$ws.Cells.Item(94, 8).Value = 1499.5  # H94
$ws.Cells.Item(94, 9).Value = 750.5  # I94
$ws.Cells.Item(94, 10).Value = 1686.75  # J94
$ws.Cells.Item(94, 11).Value = 750.5  # K94
$ws.Cells.Item(94, 12).Value = 1686.75  # L94
$ws.Cells.Item(94, 13).Value = -299.5  # M94
$ws.Cells.Item(94, 14).Value = -2588.75  # N94

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 25000  # H68
$ws.Cells.Item(68, 10).Value = 25000  # J68
$ws.Cells.Item(68, 12).Value = 25000  # L68
$ws.Cells.Item(68, 14).Value = -26498  # N68
$ws.Cells.Item(71, 8).Value = 25000  # H71
$ws.Cells.Item(71, 10).Value = 25000  # J71
$ws.Cells.Item(71, 12).Value = 75000  # L71
$ws.Cells.Item(71, 14).Value = -82488  # N71
$ws.Cells.Item(74, 8).Value = 23817.625  # H74
$ws.Cells.Item(74, 9).Value = 9000  # I74
$ws.Cells.Item(74, 10).Value = 24805.467  # J74
$ws.Cells.Item(74, 11).Value = 9000  # K74
$ws.Cells.Item(74, 12).Value = 24805.467  # L74
$ws.Cells.Item(74, 13).Value = -8126  # M74
$ws.Cells.Item(74, 14).Value = -26553.467  # N74
$ws.Cells.Item(77, 8).Value = 23817.625  # H77
$ws.Cells.Item(77, 9).Value = 9000  # I77
$ws.Cells.Item(77, 10).Value = 24805.467  # J77
$ws.Cells.Item(77, 11).Value = 27000  # K77
$ws.Cells.Item(77, 12).Value = 74416.401  # L77
$ws.Cells.Item(77, 13).Value = -22632  # M77
$ws.Cells.Item(77, 14).Value = -83152.401  # N77

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 709.1429000000001  # H113
$ws.Cells.Item(113, 9).Value = 521.3333  # I113
$ws.Cells.Item(113, 10).Value = 850  # J113
$ws.Cells.Item(113, 11).Value = 1563.9999  # K113
$ws.Cells.Item(113, 12).Value = 2550  # L113
$ws.Cells.Item(113, 13).Value = 606.0001  # M113
$ws.Cells.Item(113, 14).Value = -6890  # N113
$ws.Cells.Item(131, 8).Value = 1426.909  # H131
$ws.Cells.Item(131, 10).Value = 1852.3572  # J131
$ws.Cells.Item(131, 12).Value = 5557.071599999999  # L131
$ws.Cells.Item(131, 14).Value = -15637.0716  # N131
$ws.Cells.Item(133, 8).Value = 2856.25  # H133
$ws.Cells.Item(133, 9).Value = 2600  # I133
$ws.Cells.Item(133, 10).Value = 6700  # J133
$ws.Cells.Item(133, 11).Value = 7800  # K133
$ws.Cells.Item(133, 12).Value = 20100  # L133
$ws.Cells.Item(133, 13).Value = -2740  # M133
$ws.Cells.Item(133, 14).Value = -30220  # N133
$ws.Cells.Item(134, 8).Value = 3591.3845  # H134
$ws.Cells.Item(134, 9).Value = 1148.5  # I134
$ws.Cells.Item(134, 10).Value = 7500  # J134
$ws.Cells.Item(134, 11).Value = 3445.5  # K134
$ws.Cells.Item(134, 12).Value = 22500  # L134
$ws.Cells.Item(134, 13).Value = 1624.5  # M134
$ws.Cells.Item(134, 14).Value = -32640  # N134

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2908.1555  # H80
$ws.Cells.Item(80, 9).Value = 2442.7144  # I80
$ws.Cells.Item(80, 10).Value = 3315.4167  # J80
$ws.Cells.Item(80, 11).Value = 2442.7144  # K80
$ws.Cells.Item(80, 12).Value = 3315.4167  # L80
$ws.Cells.Item(80, 13).Value = -1444.7144  # M80
$ws.Cells.Item(80, 14).Value = -5311.4167  # N80
$ws.Cells.Item(83, 8).Value = 2908.1555  # H83
$ws.Cells.Item(83, 9).Value = 2442.7144  # I83
$ws.Cells.Item(83, 10).Value = 3315.4167  # J83
$ws.Cells.Item(83, 11).Value = 12213.572  # K83
$ws.Cells.Item(83, 12).Value = 16577.0835  # L83
$ws.Cells.Item(83, 13).Value = -7221.572  # M83
$ws.Cells.Item(83, 14).Value = -26561.0835  # N83
$ws.Cells.Item(107, 8).Value = 908.1539  # H107
$ws.Cells.Item(107, 9).Value = 633.6667  # I107
$ws.Cells.Item(107, 10).Value = 1525.75  # J107
$ws.Cells.Item(107, 11).Value = 633.6667  # K107
$ws.Cells.Item(107, 12).Value = 1525.75  # L107
$ws.Cells.Item(107, 13).Value = 1286.3333  # M107
$ws.Cells.Item(107, 14).Value = -5365.75  # N107
$ws.Cells.Item(132, 8).Value = 59737.17  # H132
$ws.Cells.Item(132, 9).Value = 68635.2  # I132
$ws.Cells.Item(132, 10).Value = 53063.65  # J132
$ws.Cells.Item(132, 11).Value = 205905.6  # K132
$ws.Cells.Item(132, 12).Value = 159190.95  # L132
$ws.Cells.Item(132, 13).Value = -203375.6  # M132
$ws.Cells.Item(132, 14).Value = -164250.95  # N132

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2095.6924  # H7
$ws.Cells.Item(7, 9).Value = 2003.5454  # I7
$ws.Cells.Item(7, 10).Value = 2602.5  # J7
$ws.Cells.Item(7, 11).Value = 2003.5454  # K7
$ws.Cells.Item(7, 12).Value = 2602.5  # L7
$ws.Cells.Item(7, 13).Value = -1891.5454  # M7
$ws.Cells.Item(7, 14).Value = -2826.5  # N7
$ws.Cells.Item(55, 8).Value = 163.27272  # H55
$ws.Cells.Item(55, 9).Value = 57.5  # I55
$ws.Cells.Item(55, 10).Value = 290.2  # J55
$ws.Cells.Item(55, 11).Value = 57.5  # K55
$ws.Cells.Item(55, 12).Value = 290.2  # L55
$ws.Cells.Item(55, 13).Value = 115.5  # M55
$ws.Cells.Item(55, 14).Value = -636.2  # N55
$ws.Cells.Item(61, 8).Value = 1916.9166  # H61
$ws.Cells.Item(61, 9).Value = 2011.5555  # I61
$ws.Cells.Item(61, 10).Value = 1633  # J61
$ws.Cells.Item(61, 11).Value = 2011.5555  # K61
$ws.Cells.Item(61, 12).Value = 1633  # L61
$ws.Cells.Item(61, 13).Value = -1809.5555  # M61
$ws.Cells.Item(61, 14).Value = -2037  # N61
$ws.Cells.Item(100, 8).Value = 30723.324  # H100
$ws.Cells.Item(100, 9).Value = 38170.11  # I100
$ws.Cells.Item(100, 11).Value = 38170.11  # K100
$ws.Cells.Item(100, 13).Value = -37629.11  # M100
$ws.Cells.Item(113, 8).Value = 1916.9166  # H113
$ws.Cells.Item(113, 9).Value = 2011.5555  # I113
$ws.Cells.Item(113, 10).Value = 1633  # J113
$ws.Cells.Item(113, 11).Value = 2011.5555  # K113
$ws.Cells.Item(113, 12).Value = 1633  # L113
$ws.Cells.Item(113, 13).Value = 158.4445000000001  # M113
$ws.Cells.Item(113, 14).Value = -5973  # N113
$ws.Cells.Item(122, 8).Value = 3120.4  # H122
$ws.Cells.Item(122, 9).Value = 3120.4  # I122
$ws.Cells.Item(122, 11).Value = 9361.200000000001  # K122
$ws.Cells.Item(122, 13).Value = -6911.200000000001  # M122
$ws.Cells.Item(126, 8).Value = 2095.6924  # H126
$ws.Cells.Item(126, 9).Value = 2003.5454  # I126
$ws.Cells.Item(126, 10).Value = 2602.5  # J126
$ws.Cells.Item(126, 11).Value = 6010.6362  # K126
$ws.Cells.Item(126, 12).Value = 7807.5  # L126
$ws.Cells.Item(126, 13).Value = -3540.6362  # M126
$ws.Cells.Item(126, 14).Value = -12747.5  # N126
$ws.Cells.Item(132, 8).Value = 61910.35  # H132
$ws.Cells.Item(132, 9).Value = 2479.8  # I132
$ws.Cells.Item(132, 10).Value = 86673.086  # J132
$ws.Cells.Item(132, 11).Value = 7439.400000000001  # K132
$ws.Cells.Item(132, 12).Value = 260019.258  # L132
$ws.Cells.Item(132, 13).Value = -4909.400000000001  # M132
$ws.Cells.Item(132, 14).Value = -265079.258  # N132

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 6874.75  # H41
$ws.Cells.Item(41, 10).Value = 6874.75  # J41
$ws.Cells.Item(41, 12).Value = 6874.75  # L41
$ws.Cells.Item(41, 14).Value = -7654.75  # N41
$ws.Cells.Item(81, 8).Value = 2772.1667  # H81
$ws.Cells.Item(81, 9).Value = 1680  # I81
$ws.Cells.Item(81, 10).Value = 2990.6  # J81
$ws.Cells.Item(81, 11).Value = 3360  # K81
$ws.Cells.Item(81, 12).Value = 5981.2  # L81
$ws.Cells.Item(81, 13).Value = -2299  # M81
$ws.Cells.Item(81, 14).Value = -8103.2  # N81
$ws.Cells.Item(84, 8).Value = 2772.1667  # H84
$ws.Cells.Item(84, 9).Value = 1680  # I84
$ws.Cells.Item(84, 10).Value = 2990.6  # J84
$ws.Cells.Item(84, 11).Value = 16800  # K84
$ws.Cells.Item(84, 12).Value = 29906  # L84
$ws.Cells.Item(84, 13).Value = -11496  # M84
$ws.Cells.Item(84, 14).Value = -40514  # N84
$ws.Cells.Item(132, 8).Value = 72938.32000000001  # H132
$ws.Cells.Item(132, 9).Value = 53448.156  # I132
$ws.Cells.Item(132, 10).Value = 114084.22  # J132
$ws.Cells.Item(132, 11).Value = 160344.468  # K132
$ws.Cells.Item(132, 12).Value = 342252.66  # L132
$ws.Cells.Item(132, 13).Value = -157814.468  # M132
$ws.Cells.Item(132, 14).Value = -347312.66  # N132
